$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38: replace "Loss to follow up" with "Not scanned" and updated description
# Row 38
$ws.Range("A38").Value = "calc_ldct_count_groups"
$ws.Range("B38").Value = "Not scanned"
$ws.Range("C38").Value = "NA"
$ws.Range("D38").Value = "The number of participants who were not scanned as part of the TLHC programme"
$ws.Range("E38").Value = "The number of unique participants in the invites table who have no follow-up information recorded in either the LHC, Measurements or LDCT tables."
$ws.Rows.Item(38).RowHeight = 30

# Row 39
$ws.Range("A39").Value = "cancer_outcome"
$ws.Range("B39").Value = "Scanned: No lung cancer"
$ws.Range("C39").Value = "LDCT AND NCRAS consolidated cancer outcomes"
$ws.Range("D39").Value = "The number of participants who had at least one low-dose CT scan and do not appear in the NCRAS dataset with a diagnosis of lung cancer."
$ws.Range("E39").Value = "The number of unique participants in the LDCT table with at least one valid scan records (dated and outcomed as being performed) but do not appear in the NCRAS dataset with a lung cancer diagnosis."
$ws.Rows.Item(39).RowHeight = 45

# Row 40
$ws.Range("A40").Value = "cancer_outcome"
$ws.Range("B40").Value = "TLHC: lung cancer"
$ws.Range("C40").Value = "LHC AND LDCT AND NCRAS consolidated cancer outcomes"
$ws.Range("D40").Value = "The number of participants with a low-dose CT scan or were assessed as high risk at LHC and a lung cancer diagnosed within 147 days of their TLHC contact."
$ws.Range("E40").Value = "The number of participants who had either a) a low-dose CT scan or b) a LHC at which they were assessed as being high risk and eligible for a scan, and also have a lung cancer diagnosed within 147 days of their TLHC contact  in the NCRAS dataset."
$ws.Rows.Item(40).RowHeight = 60

# Row 41
$ws.Range("A41").Value = "cancer_outcome"
$ws.Range("B41").Value = "No lung cancer"
$ws.Range("C41").Value = "NCRAS consolidated cancer outcomes"
$ws.Range("D41").Value = "The number of participants who were not scanned as part of the TLHC programme and do not have a lung cancer diagnosis."
$ws.Range("E41").Value = "The number of unique participants in the invites table who have no follow-up information recorded in either the LHC, Measurements or LDCT tables and are not found with a lung cancer diagnosis in the NCRAS dataset."
$ws.Rows.Item(41).RowHeight = 45

# Row 42
$ws.Range("A42").Value = "cancer_outcome"
$ws.Range("B42").Value = "Counterfactual: lung cancer"
$ws.Range("C42").Value = "Invites AND LHC AND Measurements AND LDCT AND NCRAS consolidated cancer outcomes"
$ws.Range("D42").Value = "The number of people who have a lung cancer diagnosis which is not associated with TLHC activity because they did not take up the offer of a LHC, or attended LHC but were assessed as low risk, or even had a scan but the diagnosis was made over 147 days following their scan."
$ws.Range("E42").Value = "The number of unique particiapnts in the invites table who have a lung cancer diagnosis in the NCRAS dataset but which is not associated with TLHC activity either because the participant didn't receive a scan (were invited but didn't take up the offer, attended LHC but were assessed as low risk) or did receive a scan but the diagnosis was made over 147 days afterwards."
$ws.Rows.Item(42).RowHeight = 105

# Row 43
$ws.Range("A43").Value = "cancer_stage"
$ws.Range("B43").Value = "TLHC: S 1-2"
$ws.Range("C43").Value = "NCRAS consolidated cancer outcomes"
$ws.Range("D43").Value = "The number of people with a TLHC-associated lung cancer, with a stage of either 1 or 2."
$ws.Range("E43").Value = "The number of unique participants with a lung cancer diagnosis where the diagnosis was made within 147 days of a TLHC scan or LHC at which they were assessed as high risk and with a Tumour-Node-Metastasis (TNM) staging of either 1 or 2."
$ws.Rows.Item(43).RowHeight = 60

# Row 44
$ws.Range("A44").Value = "cancer_stage"
$ws.Range("B44").Value = "TLHC: S 3-4"
$ws.Range("C44").Value = "NCRAS consolidated cancer outcomes"
$ws.Range("D44").Value = "The number of people with a TLHC-associated lung cancer, with a stage of either 3 or 4."
$ws.Range("E44").Value = "The number of unique participants with a lung cancer diagnosis where the diagnosis was made within 147 days of a TLHC scan or LHC at which they were assessed as high risk and with a Tumour-Node-Metastasis (TNM) staging of either 3 or 4."
$ws.Rows.Item(44).RowHeight = 60

# Row 45
$ws.Range("A45").Value = "cancer_stage"
$ws.Range("B45").Value = "TLHC: S ?"
$ws.Range("C45").Value = "NCRAS consolidated cancer outcomes"
$ws.Range("D45").Value = "The number of people with a TLHC-associated lung cancer which is not staged."
$ws.Range("E45").Value = "The number of unique participants with a lung cancer diagnosis where the diagnosis was made within 147 days of a TLHC scan or LHC at which they were assessed as high risk and staging information is not provided because there is insufficient information or the cancer is unstageable."
$ws.Rows.Item(45).RowHeight = 60

# Row 46
$ws.Range("A46").Value = "cancer_stage"
$ws.Range("B46").Value = "C: S 1-2"
$ws.Range("C46").Value = "NCRAS consolidated cancer outcomes"
$ws.Range("D46").Value = "The number of people with a lung cancer staged at 1 or 2 and is not associated with TLHC activity (counterfactual)."
$ws.Range("E46").Value = "The number of unique participants with a lung cancer diagnosis with a Tumour-Node-Metastasis (TNM) staging of either 1 or 2 and where the diagnosis is not associated with TLHC activity because the participant did not undergo a scan or the diagnosis was made over 147 days following scan."
$ws.Rows.Item(46).RowHeight = 60

# Row 47
$ws.Range("A47").Value = "cancer_stage"
$ws.Range("B47").Value = "C: S 3-4"
$ws.Range("C47").Value = "NCRAS consolidated cancer outcomes"
$ws.Range("D47").Value = "The number of people with a lung cancer staged at 3 or 4 and is not associated with TLHC activity (counterfactual)."
$ws.Range("E47").Value = "The number of unique participants with a lung cancer diagnosis with a Tumour-Node-Metastasis (TNM) staging of either 3 or 4 and where the diagnosis is not associated with TLHC activity because the participant did not undergo a scan or the diagnosis was made over 147 days following scan."
$ws.Rows.Item(47).RowHeight = 60

# Row 48
$ws.Range("A48").Value = "cancer_stage"
$ws.Range("B48").Value = "C: S ?"
$ws.Range("C48").Value = "NCRAS consolidated cancer outcomes"
$ws.Range("D48").Value = "The number of people with a lung cancer which is not staged and is not associated with TLHC activity (counterfactual)."
$ws.Range("E48").Value = "The number of unique participants with a lung cancer diagnosis where the diagnosis is without staging information because there is insufficient information or the cancer is unstageable and the diagnosis is not associated with TLHC activity because the participant did not undergo a scan or the diagnosis was made over 147 days following scan."
$ws.Rows.Item(48).RowHeight = 75

# Expand the Excel Table (ListObject) to cover the new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E48"))

# Restore final selection state
$ws.Range("B39").Select()
